# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# The "K" column (column G) values were recomputed and need to be
# rewritten with their new values for every data row (rows 2-29).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new value for column G ("K")
$newK = @{
    2  = 0
    3  = 3
    4  = 3
    5  = 2
    6  = 1
    7  = 5
    8  = 3
    9  = 3
    10 = 5
    11 = 1
    12 = 1
    13 = 8
    14 = 3
    15 = 4
    16 = 3
    17 = 5
    18 = 3
    19 = 2
    20 = 1
    21 = 7
    22 = 4
    23 = 5
    24 = 4
    25 = 6
    26 = 3
    27 = 3
    28 = 7
    29 = 4
}

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}
